$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $rng = $ws.Range($cell)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = $origStyle
}

Set-TextValue "D2" "34.386.60"
Set-TextValue "E2" "  +0.64%  "
Set-TextValue "D3" "1.793.40"
Set-TextValue "E3" "  +0.48%  "
Set-TextValue "E4" "  -0.16%  "
Set-TextValue "D5" "226.21"
Set-TextValue "E5" "  +0.08%  "
Set-TextValue "E6" "  +1.13%  "
Set-TextValue "E7" "  -0.15%  "
Set-TextValue "D8" "32.79"
Set-TextValue "E8" "  +2.30%  "
Set-TextValue "D9" "0.295"
Set-TextValue "E9" "  +1.22%  "
Set-TextValue "E10" "  +0.46%  "
Set-TextValue "E11" "  -0.33%  "
Set-TextValue "D12" "2.051.29"
Set-TextValue "E12" "  +0.43%  "
Set-TextValue "D13" "1.795.87"
Set-TextValue "E13" "  +0.22%  "
Set-TextValue "E14" "  +1.06%  "
Set-TextValue "E15" "  +1.59%  "
Set-TextValue "D16" "34.374.28"
Set-TextValue "E16" "  +0.61%  "
Set-TextValue "D17" "4.29"
Set-TextValue "E17" "  +2.58%  "
Set-TextValue "D18" "68.40"
Set-TextValue "E18" "  +1.08%  "
Set-TextValue "D19" "0.0₃0795"
Set-TextValue "E19" "  +0.16%  "
Set-TextValue "D20" "244.53"
Set-TextValue "E20" "  -0.45%  "
Set-TextValue "D21" "11.27"
Set-TextValue "E21" "  +2.69%  "
Set-TextValue "E22" "  -0.08%  "
Set-TextValue "E23" "  +0.77%  "
Set-TextValue "B24" "Monero"
Set-TextValue "C24" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D24" "166.59"
Set-TextValue "E24" "  +3.17%  "
Set-TextValue "B25" "Toncoin"
Set-TextValue "C25" "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue "D25" "2.07"
Set-TextValue "E25" "  +1.27%  "
Set-TextValue "D26" "7.31"
Set-TextValue "E26" "  +2.24%  "
Set-TextValue "D27" "16.51"
Set-TextValue "E27" "  +1.31%  "
Set-TextValue "E28" "  +1.17%  "
Set-TextValue "E29" "  -0.26%  "
Set-TextValue "D30" "3.99"
Set-TextValue "E30" "  +6.74%  "
Set-TextValue "D31" "0.0526"
Set-TextValue "E31" "  +1.33%  "
Set-TextValue "E32" "  +1.93%  "
Set-TextValue "E33" "  -0.01%  "
Set-TextValue "D34" "1.82"
Set-TextValue "E34" "  +1.40%  "
Set-TextValue "E35" "  -0.76%  "
Set-TextValue "D36" "1.402.52"
Set-TextValue "E36" "  -2.94%  "
Set-TextValue "D37" "0.674"
Set-TextValue "E37" "  +3.11%  "
Set-TextValue "E38" "  +2.54%  "
Set-TextValue "E39" "  -0.38%  "
Set-TextValue "D40" "85.45"
Set-TextValue "E40" "  +3.92%  "
Set-TextValue "E41" "  +4.27%  "
Set-TextValue "E42" "  +1.01%  "
Set-TextValue "E43" "  +2.40%  "
Set-TextValue "D44" "13.81"
Set-TextValue "E44" "  -2.56%  "
Set-TextValue "E45" "  +1.81%  "
Set-TextValue "E46" "  +2.93%  "
Set-TextValue "D47" "6.03"
Set-TextValue "E47" "  -0.90%  "
Set-TextValue "D48" "1.951.12"
Set-TextValue "E48" "  +0.44%  "
Set-TextValue "D49" "105.03"
Set-TextValue "E49" "  +0.20%  "
Set-TextValue "E50" "  -0.17%  "
Set-TextValue "E51" "  -2.02%  "
